$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "2025/10/08"
$ws.Range("A78").ClearFormats()
$ws.Range("B78").Value = "水"
$ws.Range("C78").Value = 10
$ws.Range("D78").Value = 201
